# ScriptOverzicht Ingepland + Charachter Controller Update
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New legend column (M3:M5): names of people scheduled ---
$ws.Range("M3").Value = "Sven"
$ws.Range("M4").Value = "Danial"
$ws.Range("M5").Value = "Harrold"

# --- New column M width ---
$ws.Range("M1").ColumnWidth = 17.71

# Colour codes used in the planning legend / table:
#   Blue   (00B0F0) -> Sven
#   Orange (FFC000)  -> Danial
#   Yellow (FFFF00)  -> Harrold
#   White  (theme Background 1) -> header cell A3

$blueColor   = 15773696   # RGB(00,B0,F0)
$orangeColor = 49407      # RGB(FF,C0,00)
$yellowColor = 65535      # RGB(FF,FF,00)

# Header cell A3 gets a plain white (theme background) fill
$ws.Range("A3").Interior.ThemeColor = 2

# Legend swatch cells
$ws.Range("M3").Interior.Color = $blueColor
$ws.Range("M4").Interior.Color = $orangeColor
$ws.Range("M5").Interior.Color = $yellowColor

# Blue (Sven) cells
$ws.Range("A4,C4,J4,A5,F5,H5,J5,A6,B6,F6,J6,A7,B7,F7,H7,J7,A8,J8,A9,J9,A10,A11,A12,A13").Interior.Color = $blueColor

# Orange (Danial) cells
$ws.Range("F4,G4,H4,C5,D5,E5,G5,C6,E6,G6,H6,D7,E7,D8,G8,D9,G9,G10").Interior.Color = $orangeColor

# Yellow (Harrold) cells
$ws.Range("B4,D4,E4,I4,B5,I5,D6,I6,G7,I7,E8,H8,I8,E9,I9,D10,I10").Interior.Color = $yellowColor

# --- View state: scroll so column B is left-most visible, selection on J9 ---
$ws.Range("J9").Select()
$excel.ActiveWindow.ScrollColumn = 2
